$wb = $excel.ActiveWorkbook

# Fix heat rate modeling syntax: update recalculated dispatch/cost/capacity values
# across the affected output sheets to match the corrected model run.

# --- Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Cells.Item(2, 2).Value = 76155.74719999997
$ws.Cells.Item(2, 3).Value = 29800
$ws.Cells.Item(2, 4).Value = 9285.87246355678
$ws.Cells.Item(2, 5).Value = 2365
$ws.Cells.Item(2, 6).Value = 15059.36469884649

# --- Capacities ---
$ws = $wb.Worksheets.Item("Capacities")
$ws.Cells.Item(3, 3).Value = 104
$ws.Cells.Item(4, 2).Value = 149
$ws.Cells.Item(4, 3).Value = 161

# --- PV Dispatch ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Cells.Item(2, 7).Value = 20.8
$ws.Cells.Item(2, 8).Value = 41.6
$ws.Cells.Item(2, 9).Value = 47.20404040404043
$ws.Cells.Item(2, 10).Value = 62.4
$ws.Cells.Item(2, 11).Value = 72.8
$ws.Cells.Item(2, 12).Value = 83.2
$ws.Cells.Item(2, 13).Value = 93.59999999999999
$ws.Cells.Item(2, 14).Value = 104
$ws.Cells.Item(2, 15).Value = 93.59999999999999
$ws.Cells.Item(2, 16).Value = 83.2
$ws.Cells.Item(2, 17).Value = 72.8
$ws.Cells.Item(2, 18).Value = 52
$ws.Cells.Item(2, 19).Value = 31.2
$ws.Cells.Item(2, 20).Value = 20.8
$ws.Cells.Item(3, 10).Value = 62.4
$ws.Cells.Item(3, 11).Value = 83.2
$ws.Cells.Item(3, 13).Value = 83.85311702887461
$ws.Cells.Item(3, 14).Value = 26
$ws.Cells.Item(3, 15).Value = 72.8
$ws.Cells.Item(3, 16).Value = 52
$ws.Cells.Item(3, 17).Value = 52
$ws.Cells.Item(3, 18).Value = 31.2
$ws.Cells.Item(3, 19).Value = 20.8
$ws.Cells.Item(4, 10).Value = 10.4
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 83.2
$ws.Cells.Item(4, 14).Value = 29.58312417100293
$ws.Cells.Item(4, 15).Value = 72.8
$ws.Cells.Item(4, 16).Value = 41.6
$ws.Cells.Item(4, 17).Value = 20.8
$ws.Cells.Item(4, 18).Value = 10.4

# --- Battery Input ---
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Cells.Item(2, 7).Value = 13
$ws.Cells.Item(2, 8).Value = 28.6
$ws.Cells.Item(2, 9).Value = 16.00404040404043
$ws.Cells.Item(2, 10).Value = 23.4
$ws.Cells.Item(2, 11).Value = 46.8
$ws.Cells.Item(2, 12).Value = 62.4
$ws.Cells.Item(2, 13).Value = 70.2
$ws.Cells.Item(2, 14).Value = 78
$ws.Cells.Item(2, 15).Value = 62.4
$ws.Cells.Item(2, 16).Value = 54.6
$ws.Cells.Item(2, 17).Value = 46.8
$ws.Cells.Item(2, 18).Value = 18.2
$ws.Cells.Item(3, 10).Value = 62.4
$ws.Cells.Item(3, 11).Value = 83.2
$ws.Cells.Item(3, 13).Value = 60.45311702887461
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 72.8
$ws.Cells.Item(3, 16).Value = 23.4
$ws.Cells.Item(3, 17).Value = 26
$ws.Cells.Item(3, 18).Value = 31.2
$ws.Cells.Item(4, 10).Value = 10.4
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 59.8
$ws.Cells.Item(4, 14).Value = 29.58312417100293
$ws.Cells.Item(4, 15).Value = 72.8
$ws.Cells.Item(4, 16).Value = 41.6
$ws.Cells.Item(4, 17).Value = 20.8
$ws.Cells.Item(4, 18).Value = 10.4

# --- Battery Output ---
$ws = $wb.Worksheets.Item("Battery Output")
$ws.Cells.Item(2, 19).Value = 10.4
$ws.Cells.Item(2, 20).Value = 25.14800000000005
$ws.Cells.Item(3, 19).Value = 20.8

# --- State of Charge ---
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Cells.Item(2, 2).Value = 187.8909090909091
$ws.Cells.Item(2, 3).Value = 168.1939393939394
$ws.Cells.Item(2, 4).Value = 155.0626262626263
$ws.Cells.Item(2, 5).Value = 141.9313131313131
$ws.Cells.Item(2, 6).Value = 128.8
$ws.Cells.Item(2, 7).Value = 141.67
$ws.Cells.Item(2, 8).Value = 169.984
$ws.Cells.Item(2, 9).Value = 185.828
$ws.Cells.Item(2, 10).Value = 208.994
$ws.Cells.Item(2, 11).Value = 255.326
$ws.Cells.Item(2, 12).Value = 317.102
$ws.Cells.Item(2, 13).Value = 386.6
$ws.Cells.Item(2, 14).Value = 463.82
$ws.Cells.Item(2, 15).Value = 525.596
$ws.Cells.Item(2, 16).Value = 579.65
$ws.Cells.Item(2, 17).Value = 625.982
$ws.Cells.Item(2, 18).Value = 644
$ws.Cells.Item(2, 19).Value = 633.4949494949495
$ws.Cells.Item(2, 20).Value = 608.0929292929293
$ws.Cells.Item(2, 21).Value = 489.9111111111111
$ws.Cells.Item(2, 22).Value = 391.4262626262627
$ws.Cells.Item(2, 23).Value = 312.6383838383838
$ws.Cells.Item(2, 24).Value = 260.1131313131313
$ws.Cells.Item(2, 25).Value = 220.7191919191919
$ws.Cells.Item(3, 2).Value = 181.3252525252525
$ws.Cells.Item(3, 3).Value = 161.6282828282828
$ws.Cells.Item(3, 4).Value = 148.4969696969697
$ws.Cells.Item(3, 5).Value = 148.4969696969697
$ws.Cells.Item(3, 6).Value = 148.4969696969697
$ws.Cells.Item(3, 7).Value = 128.8
$ws.Cells.Item(3, 8).Value = 128.8
$ws.Cells.Item(3, 9).Value = 128.8
$ws.Cells.Item(3, 10).Value = 190.576
$ws.Cells.Item(3, 11).Value = 272.944
$ws.Cells.Item(3, 12).Value = 272.944
$ws.Cells.Item(3, 13).Value = 332.7925858585859
$ws.Cells.Item(3, 14).Value = 332.7925858585859
$ws.Cells.Item(3, 15).Value = 404.8645858585859
$ws.Cells.Item(3, 16).Value = 428.0305858585859
$ws.Cells.Item(3, 17).Value = 453.7705858585859
$ws.Cells.Item(3, 18).Value = 484.6585858585859
$ws.Cells.Item(3, 19).Value = 463.6484848484848
$ws.Cells.Item(3, 20).Value = 332.3353535353535
$ws.Cells.Item(3, 21).Value = 332.3353535353535
$ws.Cells.Item(3, 22).Value = 332.3353535353535
$ws.Cells.Item(3, 23).Value = 253.5474747474748
$ws.Cells.Item(3, 24).Value = 253.5474747474748
$ws.Cells.Item(3, 25).Value = 214.1535353535353
$ws.Cells.Item(4, 2).Value = 168.1939393939394
$ws.Cells.Item(4, 3).Value = 148.4969696969697
$ws.Cells.Item(4, 4).Value = 148.4969696969697
$ws.Cells.Item(4, 5).Value = 148.4969696969697
$ws.Cells.Item(4, 6).Value = 148.4969696969697
$ws.Cells.Item(4, 7).Value = 128.8
$ws.Cells.Item(4, 8).Value = 128.8
$ws.Cells.Item(4, 9).Value = 128.8
$ws.Cells.Item(4, 10).Value = 139.096
$ws.Cells.Item(4, 11).Value = 139.096
$ws.Cells.Item(4, 12).Value = 139.096
$ws.Cells.Item(4, 13).Value = 198.298
$ws.Cells.Item(4, 14).Value = 227.5852929292929
$ws.Cells.Item(4, 15).Value = 299.6572929292929
$ws.Cells.Item(4, 16).Value = 340.8412929292929
$ws.Cells.Item(4, 17).Value = 361.4332929292929
$ws.Cells.Item(4, 18).Value = 371.7292929292929
$ws.Cells.Item(4, 19).Value = 371.7292929292929
$ws.Cells.Item(4, 20).Value = 240.4161616161616
$ws.Cells.Item(4, 21).Value = 240.4161616161616
$ws.Cells.Item(4, 22).Value = 240.4161616161616
$ws.Cells.Item(4, 23).Value = 240.4161616161616
$ws.Cells.Item(4, 24).Value = 240.4161616161616
$ws.Cells.Item(4, 25).Value = 201.0222222222222

# --- DG Dispatch ---
$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Cells.Item(2, 20).Value = 6.05199999999995
